$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting Intent/Locale right
$ws.Columns.Item(2).Insert()

# New header for inserted column
$ws.Cells.Item(1, 2).Value = "Category"

# Match formatting of the other header cells (bold font, grey fill, thin border, centered)
$ws.Cells.Item(1, 2).Font.Bold = $true
$ws.Cells.Item(1, 2).Interior.Color = 13421772
$ws.Cells.Item(1, 2).Borders.LineStyle = 1
$ws.Cells.Item(1, 2).HorizontalAlignment = -4108

# Select row 2 (where the next data row would go) to match post-edit selection
$ws.Range("A2:XFD2").Select()
